$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
# A3: 590 -> 290 (numeric-looking text; apostrophe keeps it text, ClearFormats
#     strips the auto "Text" number-format style Excel stamps on so the cell
#     style index is unaffected)
$ws.Range("A3").Value = "'290"
$ws.Range("A3").ClearFormats()

# B3: wódka Soplica -> wódka Lubelska (plain text, no numeric coercion risk)
$ws.Range("B3").Value = "wódka Lubelska"

# C3: 400 -> 800
$ws.Range("C3").Value = "'800"
$ws.Range("C3").ClearFormats()

# D3: 200 -> 300
$ws.Range("D3").Value = "'300"
$ws.Range("D3").ClearFormats()

# --- Row 4 ---
# A4: 591 -> 301
$ws.Range("A4").Value = "'301"
$ws.Range("A4").ClearFormats()

# B4: wino Fresco -> piwo Łomża
$ws.Range("B4").Value = "piwo Łomża"

# C4: 1000 -> 400
$ws.Range("C4").Value = "'400"
$ws.Range("C4").ClearFormats()

# D4: 300 -> 200
$ws.Range("D4").Value = "'200"
$ws.Range("D4").ClearFormats()

# --- Remove old rows 5 & 6 (their data was folded into rows 3 & 4 above) ---
$ws.Range("A5:E6").Delete()
